# Update crypto price/volume figures on the cryptos worksheet
# (GitHub Actions scheduled refresh).
# Values that look numeric (e.g. "7.29") are written with a leading
# apostrophe so they stay literal text cells, matching the source data
# (prices like "57.967.15" use multiple dots and are never numeric).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.967.15'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '3.123.89'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = "'529.94"
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").Value = "'138.55"
$ws.Range("E6").Value = '  -1.37%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = "'0.460"
$ws.Range("E8").Value = '  +3.35%  '
$ws.Range("D9").Value = "'7.29"
$ws.Range("E9").Value = '  +1.22%  '
$ws.Range("E10").Value = '  -1.71%  '
$ws.Range("D11").Value = "'0.406"
$ws.Range("E11").Value = '  +1.85%  '
$ws.Range("D12").Value = '3.655.60'
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("E13").Value = '  +1.19%  '
$ws.Range("D14").Value = "'25.48"
$ws.Range("E14").Value = '  -0.79%  '
$ws.Range("E15").Value = '  -2.01%  '
$ws.Range("D16").Value = '57.975.79'
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").Value = '3.117.86'
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("E18").Value = '  -2.44%  '
$ws.Range("D19").Value = "'12.65"
$ws.Range("E19").Value = '  -1.96%  '
$ws.Range("D20").Value = "'8.09"
$ws.Range("E20").Value = '  +1.30%  '
$ws.Range("D21").Value = "'351.81"
$ws.Range("E21").Value = '  -1.50%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").Value = "'68.93"
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("E24").Value = '  -1.51%  '
$ws.Range("E25").Value = '  -2.73%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").Value = '0.0₃0880'
$ws.Range("E27").Value = '  -5.71%  '
$ws.Range("D28").Value = "'7.23"
$ws.Range("E28").Value = '  -3.07%  '
$ws.Range("D29").Value = "'6.09"
$ws.Range("E29").Value = '  -4.63%  '
$ws.Range("E30").Value = '  -1.99%  '
$ws.Range("D31").Value = "'21.27"
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("D32").Value = "'4.95"
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("E33").Value = '  -4.27%  '
$ws.Range("D34").Value = "'158.86"
$ws.Range("E34").Value = '  +0.48%  '
$ws.Range("D35").Value = "'6.05"
$ws.Range("E35").Value = '  -2.43%  '
$ws.Range("D36").Value = "'26.10"
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("E37").Value = '  -2.11%  '
$ws.Range("D38").Value = "'1.66"
$ws.Range("E38").Value = '  +2.12%  '
$ws.Range("E39").Value = '  -0.65%  '
$ws.Range("E40").Value = '  -1.27%  '
$ws.Range("E41").Value = '  -3.19%  '
$ws.Range("D42").Value = "'37.48"
$ws.Range("E42").Value = '  +1.94%  '
$ws.Range("D43").Value = '2.396.80'
$ws.Range("E43").Value = '  +2.63%  '
$ws.Range("D44").Value = '3.161.96'
$ws.Range("E44").Value = '  -0.97%  '
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("E46").Value = '  -3.09%  '
$ws.Range("D47").Value = "'0.968"
$ws.Range("E47").Value = '  -2.57%  '
$ws.Range("D48").Value = "'6.03"
$ws.Range("E48").Value = '  -0.96%  '
$ws.Range("D49").Value = "'19.69"
$ws.Range("E49").Value = '  -3.52%  '
$ws.Range("D50").Value = "'0.737"
$ws.Range("E50").Value = '  -2.26%  '
$ws.Range("D51").Value = "'0.0907"
$ws.Range("E51").Value = '  +1.28%  '
